# Update countries & provincias Spain
# Applies the same edits as the source commit:
#  - Reorders a few country-name rows (the underlying data for those rows was
#    re-sorted, so both the label and that row's own numbers move together).
#  - Updates a handful of numeric case-count cells for several countries.
#  - Updates the "last updated" timestamp string.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Pais")

# ---------------------------------------------------------------------------
# 1) Update the "last updated" timestamp banner in A1
# ---------------------------------------------------------------------------
$ws.Range("A1").Value = "Datos actualizados a 12 de Abril de 2020 a las 07:22"

# ---------------------------------------------------------------------------
# 2) Country-name swaps (rows keep their position, only the label moves to
#    reflect the reordering seen in the diff; each row's own numbers are
#    updated together with the label below)
# ---------------------------------------------------------------------------

# Rows 70/71: Barein <-> Lituania
$ws.Range("A70").Value = "Lituania"
$ws.Range("A71").Value = "Barein"

# Rows 196/197/198: Nicaragua, Montserrat, Gambia -> Montserrat, Gambia, Nicaragua
$ws.Range("A196").Value = "Montserrat"
$ws.Range("A197").Value = "Gambia"
$ws.Range("A198").Value = "Nicaragua"

# Rows 200/201: Cabo Verde <-> Santa Sede
$ws.Range("A200").Value = "Santa Sede"
$ws.Range("A201").Value = "Cabo Verde"

# ---------------------------------------------------------------------------
# 3) Numeric data updates (columns B..H are: Casos totales, Nuevos casos,
#    Casos activos, Recuperados, Casos criticos, Muertes hoy, Muertes)
# ---------------------------------------------------------------------------

function Set-RowData($row, $values) {
    for ($i = 0; $i -lt $values.Length; $i++) {
        $ws.Cells.Item($row, 2 + $i).Value = $values[$i]
    }
}

Set-RowData 4   @(533088, 209, 30502, 482006, 11471, 3, 20580)   # Estados Unidos
Set-RowData 36  @(5038, 27, 1026, 3926, 37, 0, 86)                # Pakistan
Set-RowData 70  @(1053, 27, 97, 933, 14, 0, 23)                   # now Lituania
Set-RowData 71  @(1040, 0, 555, 479, 3, 0, 6)                     # now Barein
Set-RowData 75  @(880, 15, 81, 789, 21, 0, 10)                    # Kazajistan
Set-RowData 77  @(796, 29, 42, 750, 8, 0, 4)                      # Uzbekistan
Set-RowData 82  @(669, 8, 68, 573, 32, 0, 28)                     # Bulgaria
Set-RowData 92  @(501, 7, 224, 270, 17, 0, 7)                     # Uruguay
Set-RowData 196 @(9, 0, 0, 7, 0, 0, 2)                            # now Montserrat
Set-RowData 197 @(9, 0, 2, 6, 0, 0, 1)                            # now Gambia
Set-RowData 198 @(9, 0, 4, 4, 0, 0, 1)                            # now Nicaragua
Set-RowData 200 @(8, 0, 2, 6, 0, 0, 0)                            # now Santa Sede
Set-RowData 201 @(8, 0, 1, 6, 0, 0, 1)                            # now Cabo Verde
